$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (season 22/23 row already has label in A32)
$ws.Range("B32").Value = 840
$ws.Range("C32").Value = 502

$ws.Range("D32").Formula = "=B32 - C32"

# Update the view state to match: D32 becomes the active/selected cell
$ws.Range("D32").Select()
